# "Generate Report for Archive"
# Localization-status report regeneration: the handoff/status columns that
# used to read "Ready for handoff" now read "In Translation", and the
# Status/zh-cn/de-de columns are re-sized to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) / de-de (col F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value2 = $newStatus
$wsOverview.Range("F2").Value2 = $newStatus
$wsOverview.Range("E3").Value2 = $newStatus
$wsOverview.Range("F3").Value2 = $newStatus

# --- zh-cn sheet: Status column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value2 = $newStatus
$wsZhCn.Range("C3").Value2 = $newStatus

# --- de-de sheet: Status column (col C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value2 = $newStatus
$wsDeDe.Range("C3").Value2 = $newStatus

# Re-fit the status columns now that the text is shorter than
# "Ready for handoff" (the widest entry that previously drove the width).
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
